$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 730.84
$ws.Range("I15").Value = 730.84
$ws.Range("K15").Value = 2192.52
$ws.Range("M15").Value = -2023.52
$ws.Range("H53").Value = 489.5
$ws.Range("I53").Value = 472.07144
$ws.Range("J53").Value = 550.5
$ws.Range("K53").Value = 472.07144
$ws.Range("L53").Value = 550.5
$ws.Range("M53").Value = 164.92856
$ws.Range("N53").Value = -1824.5
$ws.Range("H132").Value = 8630384
$ws.Range("I132").Value = 9625759
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 28877277
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -28874747
$ws.Range("N132").Value = -16460
$ws.Range("H135").Value = 1752.3658
$ws.Range("I135").Value = 713.6667
$ws.Range("J135").Value = 2565.261
$ws.Range("K135").Value = 6423.0003
$ws.Range("L135").Value = 23087.349
$ws.Range("M135").Value = -3888.0003
$ws.Range("N135").Value = -28157.349

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2044.2354
$ws.Range("I122").Value = 2039.3334
$ws.Range("J122").Value = 2056
$ws.Range("K122").Value = 6118.0002
$ws.Range("L122").Value = 6168
$ws.Range("M122").Value = -3668.0002
$ws.Range("N122").Value = -11068
$ws.Range("H132").Value = 34063.895
$ws.Range("I132").Value = 67201.11
$ws.Range("K132").Value = 201603.33
$ws.Range("M132").Value = -199073.33

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 19558.5
$ws.Range("J35").Value = 19558.5
$ws.Range("L35").Value = 19558.5
$ws.Range("N35").Value = -20178.5
$ws.Range("H86").Value = 86735
$ws.Range("I86").Value = 93829.586
$ws.Range("K86").Value = 93829.586
$ws.Range("M86").Value = -92706.586
$ws.Range("H89").Value = 86735
$ws.Range("I89").Value = 93829.586
$ws.Range("K89").Value = 469147.93
$ws.Range("M89").Value = -463531.93

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1593
$ws.Range("I16").Value = 780
$ws.Range("J16").Value = 1999.5
$ws.Range("K16").Value = 780
$ws.Range("L16").Value = 1999.5
$ws.Range("M16").Value = -493
$ws.Range("N16").Value = -2573.5
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("H98").Value = 30266.666
$ws.Range("J98").Value = 30266.666
$ws.Range("L98").Value = 30266.666
$ws.Range("N98").Value = -34758.666
$ws.Range("H99").Value = 10827.083
$ws.Range("I99").Value = 3980
$ws.Range("J99").Value = 11449.546
$ws.Range("K99").Value = 3980
$ws.Range("L99").Value = 11449.546
$ws.Range("M99").Value = -2482
$ws.Range("N99").Value = -14445.546
$ws.Range("H113").Value = 1593
$ws.Range("I113").Value = 780
$ws.Range("J113").Value = 1999.5
$ws.Range("K113").Value = 780
$ws.Range("L113").Value = 1999.5
$ws.Range("M113").Value = 1390
$ws.Range("N113").Value = -6339.5
$ws.Range("H122").Value = 1317.7142
$ws.Range("I122").Value = 1337.3334
$ws.Range("K122").Value = 4012.0002
$ws.Range("M122").Value = -1562.0002
$ws.Range("H126").Value = 10827.083
$ws.Range("I126").Value = 3980
$ws.Range("J126").Value = 11449.546
$ws.Range("K126").Value = 11940
$ws.Range("L126").Value = 34348.638
$ws.Range("M126").Value = -9470
$ws.Range("N126").Value = -39288.638
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 125624.75
$ws.Range("J46").Value = 250749.5
$ws.Range("L46").Value = 752248.5
$ws.Range("N46").Value = -752430.5
$ws.Range("H86").Value = 720
$ws.Range("I86").Value = 450
$ws.Range("J86").Value = 900
$ws.Range("K86").Value = 1350
$ws.Range("L86").Value = 2700
$ws.Range("M86").Value = -164
$ws.Range("N86").Value = -5072
$ws.Range("H89").Value = 720
$ws.Range("I89").Value = 450
$ws.Range("J89").Value = 900
$ws.Range("K89").Value = 4050
$ws.Range("L89").Value = 8100
$ws.Range("M89").Value = 1878
$ws.Range("N89").Value = -19956
$ws.Range("H114").Value = 612.8570999999999
$ws.Range("I114").Value = 390.81818
$ws.Range("J114").Value = 1427
$ws.Range("K114").Value = 1172.45454
$ws.Range("L114").Value = 4281
$ws.Range("M114").Value = 2081.54546
$ws.Range("N114").Value = -10789

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 15000
$ws.Range("J48").Value = 15000
$ws.Range("L48").Value = 15000
$ws.Range("N48").Value = -15970
$ws.Range("H126").Value = 3924078.8
$ws.Range("I126").Value = 2494.4167
$ws.Range("J126").Value = 19610416
$ws.Range("K126").Value = 7483.250100000001
$ws.Range("L126").Value = 58831248
$ws.Range("M126").Value = -5013.250100000001
$ws.Range("N126").Value = -58836188
$ws.Range("H132").Value = 3138.923
$ws.Range("I132").Value = 2645.6
$ws.Range("J132").Value = 4783.3335
$ws.Range("K132").Value = 7936.799999999999
$ws.Range("L132").Value = 14350.0005
$ws.Range("M132").Value = -5406.799999999999
$ws.Range("N132").Value = -19410.0005

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5760.619
$ws.Range("I132").Value = 10658.333
$ws.Range("J132").Value = 3801.5334
$ws.Range("K132").Value = 31974.999
$ws.Range("L132").Value = 11404.6002
$ws.Range("M132").Value = -29444.999
$ws.Range("N132").Value = -16464.6002
$ws.Range("H136").Value = 1686.0526
$ws.Range("I136").Value = 1474.5
$ws.Range("J136").Value = 2278.4
$ws.Range("K136").Value = 4423.5
$ws.Range("L136").Value = 6835.200000000001
$ws.Range("M136").Value = -1873.5
$ws.Range("N136").Value = -11935.2

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 6805.3335
$ws.Range("J49").Value = 6937.4287
$ws.Range("L49").Value = 6937.4287
$ws.Range("N49").Value = -7397.4287
$ws.Range("H107").Value = 50565.85
$ws.Range("I107").Value = 590.86664
$ws.Range("J107").Value = 200490.8
$ws.Range("K107").Value = 1772.59992
$ws.Range("L107").Value = 601472.3999999999
$ws.Range("M107").Value = 147.4000800000001
$ws.Range("N107").Value = -605312.3999999999
$ws.Range("H108").Value = 36500
$ws.Range("J108").Value = 36500
$ws.Range("L108").Value = 36500
$ws.Range("M108").Value = -44180
$ws.Range("H113").Value = 671.6
$ws.Range("I113").Value = 499.25
$ws.Range("J113").Value = 868.5714
$ws.Range("K113").Value = 1497.75
$ws.Range("L113").Value = 2605.7142
$ws.Range("M113").Value = 672.25
$ws.Range("N113").Value = -6945.7142
$ws.Range("H122").Value = 2674.5881
$ws.Range("I122").Value = 1830.8667
$ws.Range("J122").Value = 9002.5
$ws.Range("K122").Value = 5492.6001
$ws.Range("L122").Value = 27007.5
$ws.Range("M122").Value = -3042.6001
$ws.Range("N122").Value = -31907.5
